$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 35
$ws.Cells.Item(6, 6).Value = 57
$ws.Cells.Item(7, 6).Value = 2217
$ws.Cells.Item(8, 6).Value = 63
$ws.Cells.Item(12, 6).Value = 518
$ws.Cells.Item(13, 6).Value = 1482
$ws.Cells.Item(14, 6).Value = 1482
$ws.Cells.Item(15, 6).Value = 17
$ws.Cells.Item(16, 6).Value = 574
$ws.Cells.Item(17, 6).Value = 420
$ws.Cells.Item(18, 6).Value = 807
$ws.Cells.Item(19, 6).Value = 467
$ws.Cells.Item(20, 6).Value = 3066
$ws.Cells.Item(22, 6).Value = 124
$ws.Cells.Item(23, 6).Value = 3230
$ws.Cells.Item(24, 6).Value = 688
$ws.Cells.Item(25, 6).Value = 559
$ws.Cells.Item(26, 6).Value = 257
$ws.Cells.Item(27, 6).Value = 1018
$ws.Cells.Item(28, 6).Value = 745
$ws.Cells.Item(30, 6).Value = 813
$ws.Cells.Item(31, 6).Value = 787

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 88
$ws.Cells.Item(9, 6).Value = 71
$ws.Cells.Item(20, 6).Value = 206
$ws.Cells.Item(21, 6).Value = 149

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(6, 6).Value = 427

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(8, 6).Value = 35
$ws.Cells.Item(9, 6).Value = 88
$ws.Cells.Item(12, 6).Value = 57
$ws.Cells.Item(13, 6).Value = 427
$ws.Cells.Item(14, 6).Value = 2217
$ws.Cells.Item(15, 6).Value = 63
$ws.Cells.Item(20, 6).Value = 71
$ws.Cells.Item(23, 6).Value = 518
$ws.Cells.Item(26, 6).Value = 1482
$ws.Cells.Item(27, 6).Value = 1482
$ws.Cells.Item(29, 6).Value = 17
$ws.Cells.Item(30, 6).Value = 574
$ws.Cells.Item(31, 6).Value = 420
$ws.Cells.Item(32, 6).Value = 807
$ws.Cells.Item(33, 6).Value = 467
$ws.Cells.Item(35, 6).Value = 3066
$ws.Cells.Item(36, 6).Value = 124
$ws.Cells.Item(37, 6).Value = 3230
$ws.Cells.Item(38, 6).Value = 688
$ws.Cells.Item(39, 6).Value = 559
$ws.Cells.Item(40, 6).Value = 257
$ws.Cells.Item(41, 6).Value = 1018
$ws.Cells.Item(44, 6).Value = 206
$ws.Cells.Item(45, 6).Value = 149
$ws.Cells.Item(47, 6).Value = 745
$ws.Cells.Item(49, 6).Value = 813
$ws.Cells.Item(50, 6).Value = 787
